$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2026-01-09 Friday" "2026-01-10 Saturday"

Replace-Text "785÷6=130, 5" "478÷7=68, 2"
Replace-Text "387÷9=43, 0" "186÷7=26, 4"
Replace-Text "397÷2=198, 1" "293÷5=58, 3"
Replace-Text "143÷3=47, 2" "753÷9=83, 6"
Replace-Text "418÷4=104, 2" "927÷2=463, 1"
Replace-Text "484÷8=60, 4" "896÷7=128, 0"
Replace-Text "956÷5=191, 1" "515÷8=64, 3"
Replace-Text "660÷5=132, 0" "976÷5=195, 1"
Replace-Text "675÷6=112, 3" "136÷6=22, 4"
Replace-Text "396÷3=132, 0" "676÷5=135, 1"
Replace-Text "881÷7=125, 6" "669÷4=167, 1"
Replace-Text "841÷6=140, 1" "212÷5=42, 2"
Replace-Text "424÷9=47, 1" "439÷6=73, 1"
Replace-Text "401÷5=80, 1" "510÷6=85, 0"
Replace-Text "877÷5=175, 2" "202÷9=22, 4"
Replace-Text "807÷8=100, 7" "432÷3=144, 0"
Replace-Text "823÷8=102, 7" "228÷8=28, 4"
Replace-Text "984÷2=492, 0" "563÷2=281, 1"
Replace-Text "410÷4=102, 2" "832÷8=104, 0"
Replace-Text "249÷3=83, 0" "879÷9=97, 6"
Replace-Text "799÷7=114, 1" "375÷3=125, 0"
Replace-Text "407÷7=58, 1" "437÷9=48, 5"
Replace-Text "812÷9=90, 2" "143÷6=23, 5"
Replace-Text "624÷3=208, 0" "444÷7=63, 3"
Replace-Text "174÷6=29, 0" "971÷8=121, 3"

Write-Output "Replacements complete"
